$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repulled data
$ws.Range("F2").Value = -7
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = -5
$ws.Range("F6").Value = -10
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = 0
